# Update cryptocurrency price (column D) and 1h volume change (column E) values
# as refreshed by the scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}

Set-TextValue "D2" "304.45"
Set-TextValue "E2" "5.63%"
Set-TextValue "D3" "34.90"
Set-TextValue "E3" "12.34%"
Set-TextValue "D4" "5.213"
Set-TextValue "E4" "5.85%"
Set-TextValue "D5" "0.07814"
Set-TextValue "E5" "6.61%"
Set-TextValue "D6" "2.378"
Set-TextValue "E6" "6.46%"
Set-TextValue "D7" "8.043"
Set-TextValue "E7" "4.29%"
Set-TextValue "D8" "0.9354"
Set-TextValue "E8" "3.74%"
Set-TextValue "D9" "0.1012"
Set-TextValue "E9" "10.04%"
Set-TextValue "D10" "0.1855"
Set-TextValue "E10" "9.56%"
Set-TextValue "D11" "0.08634"
Set-TextValue "E11" "4.75%"
Set-TextValue "D12" "0.03314"
Set-TextValue "E12" "6.36%"
Set-TextValue "D13" "0.09896"
Set-TextValue "E13" "-0.44%"
Set-TextValue "D14" "0.001493"
Set-TextValue "E14" "-0.34%"
Set-TextValue "D15" "0.005755"
Set-TextValue "E15" "1.09%"
Set-TextValue "E16" "-1.74%"
Set-TextValue "D17" "3.962"
Set-TextValue "E17" "6.33%"
Set-TextValue "D18" "2.148"
Set-TextValue "E18" "3.20%"
Set-TextValue "D20" "0.1302"
Set-TextValue "E20" "0.25%"
Set-TextValue "D21" "4.312"
Set-TextValue "E21" "3.89%"
Set-TextValue "D22" "0.2226"
Set-TextValue "E22" "5.92%"
Set-TextValue "D23" "0.04576"
Set-TextValue "E23" "0.89%"
Set-TextValue "D24" "0.001218"
Set-TextValue "E24" "0.66%"
Set-TextValue "D25" "0.004438"
Set-TextValue "E25" "6.59%"
Set-TextValue "E26" "-0.03%"
Set-TextValue "D27" "0.0003697"
Set-TextValue "D39" "0.01776"
Set-TextValue "E39" "12.97%"
Set-TextValue "D40" "0.04805"
Set-TextValue "E40" "8.04%"
Set-TextValue "D41" "0.007764"
Set-TextValue "E41" "4.46%"
Set-TextValue "E42" "6.11%"
Set-TextValue "D43" "0.007035"
Set-TextValue "E43" "-26.26%"
Set-TextValue "D44" "0.002299"
Set-TextValue "E44" "-0.03%"
Set-TextValue "D45" "0.009478"
Set-TextValue "E45" "4.28%"
Set-TextValue "D46" "0.00005951"
Set-TextValue "E46" "-2.68%"
Set-TextValue "D47" "0.00000000749"
Set-TextValue "E47" "-0.05%"
Set-TextValue "E48" "12.74%"
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "-0.05%"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "-0.05%"
